$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "Recovered" / "Matched on recovery?" info for rows 2-5 (I/J columns)
#    Style source: I13 (date, s=4) and J6 ("yes", s=1) already carry the
#    correct formatting, so copy formats from them first, then set values.
# ---------------------------------------------------------------------------
$ws.Range("I13").Copy()
$ws.Range("I2:I5").PasteSpecial(-4122)
$ws.Range("I25:I30").PasteSpecial(-4122)

$ws.Range("J6").Copy()
$ws.Range("J2:J5").PasteSpecial(-4122)
$ws.Range("J25:J30").PasteSpecial(-4122)

$ws.Range("I2").Value = 45250
$ws.Range("J2").Value = "yes"
$ws.Range("I3").Value = 45250
$ws.Range("J3").Value = "yes"
$ws.Range("I4").Value = 45250
$ws.Range("J4").Value = "yes"
$ws.Range("I5").Value = 45250
$ws.Range("J5").Value = "yes"
$ws.Range("I25").Value = 45250
$ws.Range("J25").Value = "yes"
$ws.Range("I26").Value = 45250
$ws.Range("J26").Value = "yes"

# ---------------------------------------------------------------------------
# 2) Re-style the existing "recovered" dates in rows 19-26 so they use the
#    10pt font (style s=4) instead of the stray 12pt one (style s=5).
# ---------------------------------------------------------------------------
$ws.Range("I19").Font.Size = 10
$ws.Range("I20").Font.Size = 10
$ws.Range("I21").Font.Size = 10
$ws.Range("I22").Font.Size = 10
$ws.Range("I23").Font.Size = 10
$ws.Range("I24").Font.Size = 10
$ws.Range("M21").Font.Size = 10
$ws.Range("M22").Font.Size = 10
$ws.Range("M23").Font.Size = 10
$ws.Range("M24").Font.Size = 10
$ws.Range("E25").Font.Size = 10
$ws.Range("E26").Font.Size = 10

# K21:L24 had no explicit style at all (12pt default); bring them to the
# same 10pt body style (s=1) used throughout the rest of the table.
$ws.Range("K21").Font.Size = 10
$ws.Range("L21").Font.Size = 10
$ws.Range("K22").Font.Size = 10
$ws.Range("L22").Font.Size = 10
$ws.Range("K23").Font.Size = 10
$ws.Range("L23").Font.Size = 10
$ws.Range("K24").Font.Size = 10
$ws.Range("L24").Font.Size = 10

# ---------------------------------------------------------------------------
# 3) Row 26 redeployment site correction: MacNorth instead of Emerald.
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = "MacNorth"

# ---------------------------------------------------------------------------
# 4) Four brand-new deployment rows (27-30), matching the layout of the
#    other "A/C/D/E" + "I/J" rows above (copy formats, then fill values).
# ---------------------------------------------------------------------------
$ws.Range("A25").Copy()
$ws.Range("A27:A30").PasteSpecial(-4122)

$ws.Range("C25").Copy()
$ws.Range("C27:C30").PasteSpecial(-4122)

$ws.Range("D25").Copy()
$ws.Range("D27:D30").PasteSpecial(-4122)

$ws.Range("E25").Copy()
$ws.Range("E27:E30").PasteSpecial(-4122)

$ws.Range("A27").Value = "CTD"
$ws.Range("C27").Value = 6635
$ws.Range("D27").Value = "Emerald"
$ws.Range("E27").Value = 45121
$ws.Range("I27").Value = 45250
$ws.Range("J27").Value = "yes"

$ws.Range("A28").Value = "CTD"
$ws.Range("C28").Value = 7398
$ws.Range("D28").Value = "MacNorth"
$ws.Range("E28").Value = 45121
$ws.Range("I28").Value = 45250
$ws.Range("J28").Value = "yes"

$ws.Range("A29").Value = "EcoPAR"
$ws.Range("C29").Value = 4298
$ws.Range("D29").Value = "Emerald"
$ws.Range("E29").Value = 45121
$ws.Range("I29").Value = 45250
$ws.Range("J29").Value = "yes"

$ws.Range("A30").Value = "EcoPAR"
$ws.Range("C30").Value = 4297
$ws.Range("D30").Value = "MacNorth"
$ws.Range("E30").Value = 45121
$ws.Range("I30").Value = 45250
$ws.Range("J30").Value = "yes"

$ws.Range("D27").Select()
